$p = $ppt.ActivePresentation

# --- Helper: replace a "Unit18" prefix (first 6 chars) with "Unit17" in a
#     shape's text range, leaving the rest of the text/runs untouched. ---
function Set-Unit18ToUnit17($shape) {
    $tr = $shape.TextFrame.TextRange
    $chars = $tr.Characters(1, 6)
    if ($chars.Text -eq "Unit18") {
        $chars.Text = "Unit17"
    }
}

# 1) Slide 1: title-slide "UNIT 18" textbox (two runs "UNIT 1" + "8") becomes
#    a single run "UNIT 17".
$s1 = $p.Slides.Item(1)
$titleShape = $s1.Shapes.Item(2)
$tr1 = $titleShape.TextFrame.TextRange
$tr1.Delete()
[void]$tr1.InsertBefore("UNIT 17")

# 2) Slide 2: big title placeholder "Unit 18: Heap" -> "Unit 17: Heap"
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Unit 17: Heap"

# 3) Slides 2-9: "Slide Number Placeholder" shape (id 7), first run "Unit18" -> "Unit17"
for ($si = 2; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.Id -eq 7 -and $sh.HasTextFrame -and $sh.TextFrame.HasText) {
            Set-Unit18ToUnit17 $sh
        }
    }
}

# 4) Slide master: "Slide Number Placeholder" shape, "Unit18" -> "Unit17"
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -like "Unit18*") {
        Set-Unit18ToUnit17 $sh
    }
}

# 5) All slide layouts: "Slide Number Placeholder" shape, "Unit18" -> "Unit17"
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -like "Unit18*") {
            Set-Unit18ToUnit17 $sh
        }
    }
}
